# Update handback status timestamps ("Generate Report for Handback")
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: Correspond Handoff/Handback DateTime for row 2
$wsZhCn.Range("E2").Value = "2016-03-23 09:19:28"
$wsZhCn.Range("H2").Value = "2016-03-23 09:19:53"

# de-de sheet: Correspond Handoff/Handback DateTime for row 2
$wsDeDe.Range("E2").Value = "2016-03-23 09:19:32"
$wsDeDe.Range("H2").Value = "2016-03-23 09:20:01"
